# update scripts wuth new tpm
#
# The NATMI ligand-receptor edge table for Dkk2->Lrp5 was recomputed with
# the refreshed TPM expression matrix. The recompute now also yields
# "ECs" as a sending cluster (in addition to the existing "FAPs" sender),
# so the table grows from 3 data rows to 6 data rows covering every
# sender/target combination, and every metric column is refreshed with
# the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:T, rows 2-7 (row 1 is the header and is unchanged).
$data = @(
    @("ECs",  "Dkk2", "Lrp5", "ECs",   1, 0.3333333333333333, 0.1293116666666667, 0.387935,  0.0151634279701504, 0.0151634279701504, 3, 1, 15.01856033333333, 45.055681, 0.4908713633047416, 0.4908713633047417, 1.942075067637222,  17.478675608735,  0.007443292560080979, 0.007443292560080981),
    @("ECs",  "Dkk2", "Lrp5", "FAPs",  1, 0.3333333333333333, 0.1293116666666667, 0.387935,  0.0151634279701504, 0.0151634279701504, 3, 1, 9.484070666666666, 28.452212, 0.3099803572711625, 0.3099803572711625, 1.226400984691111,  11.03760886222,  0.004700364819642759, 0.00470036481964276),
    @("ECs",  "Dkk2", "Lrp5", "MuSCs", 1, 0.3333333333333333, 0.1293116666666667, 0.387935,  0.0151634279701504, 0.0151634279701504, 3, 1, 6.093084,          18.279252, 0.1991482794240958, 0.1991482794240958, 0.7879068471799999, 7.09116162462,   0.003019770590426662, 0.003019770590426663),
    @("FAPs", "Dkk2", "Lrp5", "ECs",   3, 1,                  8.398553333333334, 25.19566,  0.9848365720298496, 0.9848365720298496, 3, 1, 15.01856033333333, 45.055681, 0.4908713633047416, 0.4908713633047417, 126.1341799493845,  1135.20761954446, 0.4834280707446607,   0.4834280707446607),
    @("FAPs", "Dkk2", "Lrp5", "FAPs",  3, 1,                  8.398553333333334, 25.19566,  0.9848365720298496, 0.9848365720298496, 3, 1, 9.484070666666666, 28.452212, 0.3099803572711625, 0.3099803572711625, 79.65247331110223,  716.87225979992, 0.3052799924515197,   0.3052799924515198),
    @("FAPs", "Dkk2", "Lrp5", "MuSCs", 3, 1,                  8.398553333333334, 25.19566,  0.9848365720298496, 0.9848365720298496, 3, 1, 6.093084,          18.279252, 0.1991482794240958, 0.1991482794240958, 51.17309093848001,  460.55781844632, 0.1961285088336692,   0.1961285088336692)
)

$r = 2
foreach ($rowValues in $data) {
    $c = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
